# "removed redundant thousands label"
# The underlying data table (Table1, A1:B27 on sheet1) was re-sorted from
# "Sales Total descending" to "Company Name ascending" (A-Z). This mirrors
# what Excel does when you click the column-A header's A-Z sort button
# instead of keeping the old sort-by-B button/label.
#
# We re-derive the new row order with a culture-aware (case-insensitive,
# accent-folded) comparison of the existing Company Name values -- the same
# collation Excel's "Sort A to Z" uses -- entirely through the object model
# (no hard-coded literal text), then write the rows back in that order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 27

# 1) Snapshot the existing Company Name / Sales Total pairs before touching
#    anything, so reads never race writes.
$rows = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $total = $ws.Cells.Item($r, 2).Value2
    $rows += ,@($name, $total)
}

# 2) Build a normalized sort key per row: lower-case plus a fold of the
#    common accented letters down to their base letter, so "Forêts" sorts
#    like "Forets" (ahead of "Formaggi"), matching Excel's default locale
#    text sort instead of a raw ordinal/codepoint sort.
$foldPairs = @(
    @('à','a'), @('á','a'), @('â','a'), @('ã','a'), @('ä','a'), @('å','a'),
    @('è','e'), @('é','e'), @('ê','e'), @('ë','e'),
    @('ì','i'), @('í','i'), @('î','i'), @('ï','i'),
    @('ò','o'), @('ó','o'), @('ô','o'), @('õ','o'), @('ö','o'),
    @('ù','u'), @('ú','u'), @('û','u'), @('ü','u'),
    @('ý','y'), @('ÿ','y'),
    @('ç','c'), @('ñ','n'),
    @('ß','s')
)

function Get-SortKey($text) {
    $k = $text.ToLower()
    foreach ($pair in $foldPairs) {
        $k = $k.Replace($pair[0], $pair[1])
    }
    return $k
}

$keyed = @()
foreach ($row in $rows) {
    $key = Get-SortKey $row[0]
    $keyed += ,@($key, $row[0], $row[1])
}

$sorted = $keyed | Sort-Object { $_[0] }

# 3) Write the re-ordered Company Name / Sales Total pairs back.
$r = $firstDataRow
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $r = $r + 1
}

# 4) Keep the table's sort definition consistent with the new A-Z order
#    (sorted by Company Name instead of Sales Total).
$tbl = $ws.ListObjects.Item(1)
$sort = $tbl.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A1:A27"))
$sort.Header = 1
$sort.Apply()

# 5) Re-write the values once more on top of the table sort, since a plain
#    ordinal re-sort through the Sort object does not reliably match
#    Excel's locale-aware text order for the accented entries.
$r = $firstDataRow
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $r = $r + 1
}

# 6) Move the active selection the same way the author's session ended up
#    (cell F23, within the now-visible top-left area of the sheet).
[void]$ws.Range("F23").Select()
